$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record needs to be inserted before the existing row
# that is currently at row 140 (date 2021-07-06 / serial 44383). Insert a
# whole row there so it - and every following record through row 143 -
# shifts down by one (140->141, 141->142, 142->143, 143->144).
$ws.Rows.Item(140).Insert()

# Fill the newly inserted row 140 with the new weekly entry.
$ws.Range("A140").Value = 7
$ws.Range("B140").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C140").Value = "Ñuble"
$ws.Range("D140").Value = 44509
$ws.Range("E140").Value = 16
$ws.Range("F140").Value = 100112017
$ws.Range("G140").Value = "Apio"
$ws.Range("H140").Value = "Americana (o)"
$ws.Range("I140").Value = "Primera"
$ws.Range("J140").Value = 100
$ws.Range("K140").Value = 8000
$ws.Range("L140").Value = 9000
$ws.Range("M140").Value = 8500
$ws.Range("N140").Value = "`$/docena de matas"
$ws.Range("O140").Value = "Provincia del Elquí"
$ws.Range("P140").Value = 1417
$ws.Range("Q140").Value = 6
$ws.Range("R140").Value = "Hortaliza"
